$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "890×2="; New = "319×3=" },
    @{ Old = "902×3="; New = "263×9=" },
    @{ Old = "361×3="; New = "761×7=" },
    @{ Old = "230×4="; New = "681×7=" },
    @{ Old = "472×9="; New = "674×6=" },
    @{ Old = "591×6="; New = "961×9=" },
    @{ Old = "149×8="; New = "188×7=" },
    @{ Old = "822×3="; New = "816×8=" },
    @{ Old = "609×9="; New = "435×6=" },
    @{ Old = "877×5="; New = "958×5=" },
    @{ Old = "684×8="; New = "987×6=" },
    @{ Old = "930×5="; New = "503×4=" },
    @{ Old = "955×4="; New = "290×8=" },
    @{ Old = "556×3="; New = "257×9=" },
    @{ Old = "626×2="; New = "492×8=" },
    @{ Old = "492×7="; New = "526×4=" },
    @{ Old = "182×7="; New = "127×6=" },
    @{ Old = "261×2="; New = "396×2=" },
    @{ Old = "774×6="; New = "149×5=" },
    @{ Old = "564×3="; New = "922×2=" },
    @{ Old = "635×3="; New = "697×6=" },
    @{ Old = "637×3="; New = "151×6=" },
    @{ Old = "326×8="; New = "356×4=" },
    @{ Old = "882×7="; New = "416×7=" },
    @{ Old = "644×5="; New = "797×4=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
